{"js": "// Enhance Siege Analytics \"Research & Data Analytics Leadership\" bullets with\n// three new achievement lines (voter file discovery, boundary estimation\n// algorithm, cost savings impact), inserted right after the section's lead-in\n// paragraph and before the existing \"\u2022 Conceived, architected...\" bullet.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText = \"Research & Data Analytics Leadership\";\nlet anchorParagraph = null;\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text === anchorText) {\n    anchorParagraph = paragraph;\n    break;\n  }\n}\n\nif (!anchorParagraph) {\n  throw new Error(`Could not find anchor paragraph: \"${anchorText}\"`);\n}\n\nconst newBullets = [\n  \"\u2022 Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters\",\n  \"\u2022 Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States\",\n  \"\u2022 Algorithm reduced mapping costs by 75%, saving campaigns and organizations $5M+ and enabling smaller nonprofits to conduct redistricting analysis\"\n];\n\n// Insert each bullet after the previous one so the resulting order matches\n// the source order (chaining off the freshly-inserted paragraph keeps them\n// in sequence instead of reversing).\nlet insertAfter = anchorParagraph;\nfor (const bulletText of newBullets) {\n  insertAfter = insertAfter.insertParagraph(bulletText, \"After\");\n}\n\nawait context.sync();\n", "ps1": "# Enhance Siege Analytics \"Research & Data Analytics Leadership\" bullets with\n# three new achievement lines (voter file discovery, boundary estimation\n# algorithm, cost savings impact), inserted right after the section's lead-in\n# paragraph and before the existing \"Conceived, architected...\" bullet.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"Research & Data Analytics Leadership\"\n\n$targetIndex = -1\n$i = 1\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd() -eq $anchorText) {\n        $targetIndex = $i\n        break\n    }\n    $i = $i + 1\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not find anchor paragraph: $anchorText\"\n}\n\n$newBullets = @(\n    \"\u2022 Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters\",\n    \"\u2022 Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States\",\n    \"\u2022 Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis\"\n)\n\n# Re-fetch the paragraph by its (now-shifting) index before each insert so\n# every new bullet lands immediately after the previous one, in source order.\n$insertIndex = $targetIndex\nforeach ($bulletText in $newBullets) {\n    $anchorParagraph = $d.Paragraphs.Item($insertIndex)\n    $insertRange = $anchorParagraph.Range\n    $insertRange.Collapse(0)\n    $insertRange.InsertAfter($bulletText + \"`r\")\n    $insertIndex = $insertIndex + 1\n}\n"}
